$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers to the shorter export labels
$ws.Range("F1").Value = "Memory tot.size"
$ws.Range("H1").Value = "Memory P/Ns"
$ws.Range("M1").Value = "HDD slot pop."
$ws.Range("N1").Value = "PSU P/Ns"

# Adjust column widths (geometry fix)
$ws.Columns.Item(6).ColumnWidth = 15
$ws.Columns.Item(8).ColumnWidth = 11
$ws.Columns.Item(13).ColumnWidth = 13
$ws.Columns.Item(14).ColumnWidth = 8
